# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on all three sheets
# - Refresh the associated "last generated" timestamps
# - Widen the Status-related columns so the new, longer text fits

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview "Latest HO Xliff Generate Date" (G2) and de-de "Latest Handback DateTime" (H2)
$ws_overview.Range("G2").Value = "2016-08-25 18:40:49"
$ws_dede.Range("H2").Value = "2016-08-25 18:40:49"

# zh-cn "Latest Handoff Datetime" (H2)
$ws_zhcn.Range("H2").Value = "2016-08-25 18:40:45"

# --- Column widths: widen the Status-related columns to fit "Ready for handoff" ---
$newWidth = 16.3333333333333

$ws_overview.Columns.Item(5).ColumnWidth = $newWidth   # column E
$ws_overview.Columns.Item(6).ColumnWidth = $newWidth   # column F
$ws_zhcn.Columns.Item(3).ColumnWidth = $newWidth       # column C
$ws_dede.Columns.Item(3).ColumnWidth = $newWidth       # column C
